$wb = $excel.ActiveWorkbook

# --- Metadata sheet updates ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/measure-report-evidence-population-id"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

# --- Elements sheet updates ---
$elements = $wb.Worksheets.Item("Elements")

# Fixed Value for Extension.url now reflects the updated canonical URL
$elements.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/measure-report-evidence-population-id"

# The ele-1/ext-1 constraint text now lives only on the "Extension.extension"
# row (row 4, unchanged); the base "Extension" row (row 2) no longer repeats it.
$elements.Range("AI2").Value = ""
